$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1727.3793
$ws.Range("I40").Value = 1616.3334
$ws.Range("J40").Value = 1805.7646
$ws.Range("K40").Value = 1616.3334
$ws.Range("L40").Value = 1805.7646
$ws.Range("M40").Value = -1441.3334
$ws.Range("N40").Value = -2155.7646

# Row 107
$ws.Range("H107").Value = 705.6316
$ws.Range("I107").Value = 705.6316
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 705.6316
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1214.3684
$ws.Range("N107").ClearContents()

# Row 112
$ws.Range("H112").Value = 80774.234
$ws.Range("I112").Value = 375
$ws.Range("J112").Value = 87474.164
$ws.Range("K112").Value = 1125
$ws.Range("L112").Value = 262422.492
$ws.Range("M112").Value = -17
$ws.Range("N112").Value = -264638.492

# Row 118
$ws.Range("H118").Value = 967.0909
$ws.Range("I118").Value = 967.0909
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 2901.2727
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -1244.2727

# Row 135
$ws.Range("H135").Value = 1830.8667
$ws.Range("I135").Value = 1437.7273
$ws.Range("J135").Value = 2912
$ws.Range("K135").Value = 12939.5457
$ws.Range("L135").Value = 26208
$ws.Range("M135").Value = -10404.5457
$ws.Range("N135").Value = -31278

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2640.3333
$ws.Range("I2").Value = 2237.0588
$ws.Range("J2").Value = 3167.6924
$ws.Range("K2").Value = 2237.0588
$ws.Range("L2").Value = 3167.6924
$ws.Range("M2").Value = -2124.0588
$ws.Range("N2").Value = -3393.6924

# Row 22
$ws.Range("H22").Value = 18104.857
$ws.Range("I22").Value = 3905.3333
$ws.Range("J22").Value = 28754.5
$ws.Range("K22").Value = 3905.3333
$ws.Range("L22").Value = 28754.5
$ws.Range("M22").Value = -3606.3333
$ws.Range("N22").Value = -29352.5

# Row 32
$ws.Range("H32").Value = 3411.7078
$ws.Range("I32").Value = 2966.0334
$ws.Range("J32").Value = 8759.799999999999
$ws.Range("K32").Value = 2966.0334
$ws.Range("L32").Value = 8759.799999999999
$ws.Range("M32").Value = -2679.0334

# Row 45
$ws.Range("H45").Value = 5132.3516
$ws.Range("I45").Value = 7466.9443
$ws.Range("J45").Value = 2920.6316
$ws.Range("K45").Value = 7466.9443
$ws.Range("L45").Value = 2920.6316
$ws.Range("M45").Value = -7089.9443
$ws.Range("N45").Value = -3674.6316

# Row 101
$ws.Range("H101").Value = 50000
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 50000
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 50000
$ws.Range("N101").Value = -56490

# Row 102
$ws.Range("H102").Value = 4970.5835
$ws.Range("I102").Value = 4964.7
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 4964.7
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -3342.7

# Row 116
$ws.Range("H116").Value = 2640.3333
$ws.Range("I116").Value = 2237.0588
$ws.Range("J116").Value = 3167.6924
$ws.Range("K116").Value = 2237.0588
$ws.Range("L116").Value = 3167.6924
$ws.Range("M116").Value = 56.94120000000021
$ws.Range("N116").Value = -7755.6924

# Row 132
$ws.Range("H132").Value = 3488.682
$ws.Range("I132").Value = 2782.5264
$ws.Range("J132").Value = 7961
$ws.Range("K132").Value = 8347.5792
$ws.Range("L132").Value = 23883
$ws.Range("M132").Value = -5817.5792
$ws.Range("N132").Value = -28943

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2640.3333
$ws.Range("I3").Value = 2237.0588
$ws.Range("J3").Value = 3167.6924
$ws.Range("K3").Value = 2237.0588
$ws.Range("L3").Value = 3167.6924
$ws.Range("M3").Value = -2123.0588
$ws.Range("N3").Value = -3395.6924

# Row 99
$ws.Range("H99").Value = 6091.5557
$ws.Range("I99").Value = 1764.8
$ws.Range("J99").Value = 11500
$ws.Range("K99").Value = 1764.8
$ws.Range("L99").Value = 11500
$ws.Range("M99").Value = -266.8
$ws.Range("N99").Value = -14496

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1753.6
$ws.Range("I16").Value = 1942
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 1942
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -1655

# Row 31
$ws.Range("H31").Value = 43998.32
$ws.Range("I31").Value = 57051.832
$ws.Range("J31").Value = 10432.143
$ws.Range("K31").Value = 57051.832
$ws.Range("L31").Value = 10432.143
$ws.Range("M31").Value = -56756.832
$ws.Range("N31").Value = -11022.143

# Row 34
$ws.Range("H34").Value = 43998.32
$ws.Range("I34").Value = 57051.832
$ws.Range("J34").Value = 10432.143
$ws.Range("K34").Value = 57051.832
$ws.Range("L34").Value = 10432.143
$ws.Range("M34").Value = -56849.832
$ws.Range("N34").Value = -10836.143

# Row 58
$ws.Range("H58").Value = 3017.6924
$ws.Range("I58").Value = 3113.6365
$ws.Range("J58").Value = 2490
$ws.Range("K58").Value = 3113.6365
$ws.Range("L58").Value = 2490
$ws.Range("M58").Value = -2910.6365

# Row 81
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# Row 84
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# Row 94
$ws.Range("H94").Value = 1918.875
$ws.Range("I94").Value = 1978.7142
$ws.Range("J94").Value = 1500
$ws.Range("K94").Value = 1978.7142
$ws.Range("L94").Value = 1500
$ws.Range("M94").Value = -1527.7142
$ws.Range("N94").Value = -2402

# Row 99
$ws.Range("H99").Value = 3376.125
$ws.Range("I99").Value = 3376.125
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3376.125
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1878.125
$ws.Range("N99").ClearContents()

# Row 112
$ws.Range("H112").Value = 76749.5
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 76749.5
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 76749.5
$ws.Range("N112").Value = -79703.5

# Row 113
$ws.Range("H113").Value = 1753.6
$ws.Range("I113").Value = 1942
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 1942
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 228

# Row 126
$ws.Range("H126").Value = 3376.125
$ws.Range("I126").Value = 3376.125
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 10128.375
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -7658.375
$ws.Range("N126").ClearContents()

# Row 127
$ws.Range("H127").Value = 80000
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 80000
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 80000
$ws.Range("N127").Value = -89920

# Row 132
$ws.Range("H132").Value = 3783.923
$ws.Range("I132").Value = 3381
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 10143
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -7613
$ws.Range("N132").Value = -23060

# Row 136
$ws.Range("H136").Value = 3017.6924
$ws.Range("I136").Value = 3113.6365
$ws.Range("J136").Value = 2490
$ws.Range("K136").Value = 9340.9095
$ws.Range("L136").Value = 7470
$ws.Range("M136").Value = -6790.9095

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 1067826.6
$ws.Range("I4").Value = 1167334.4
$ws.Range("J4").Value = 212060
$ws.Range("K4").Value = 3502003.2
$ws.Range("L4").Value = 636180
$ws.Range("M4").Value = -3501891.2

# Row 9
$ws.Range("H9").Value = 25994.5
$ws.Range("I9").Value = 987.25
$ws.Range("J9").Value = 38498.125
$ws.Range("K9").Value = 2961.75
$ws.Range("L9").Value = 115494.375
$ws.Range("M9").Value = -2737.75
$ws.Range("N9").Value = -115942.375

# Row 10
$ws.Range("H10").Value = 183.1
$ws.Range("I10").Value = 243.71428
$ws.Range("J10").Value = 41.666668
$ws.Range("K10").Value = 731.14284
$ws.Range("L10").Value = 125.000004
$ws.Range("M10").Value = -592.14284
$ws.Range("N10").Value = -403.000004

# Row 129
$ws.Range("H129").Value = 820.25
$ws.Range("I129").Value = 649.3333
$ws.Range("J129").Value = 1333
$ws.Range("K129").Value = 1947.9999
$ws.Range("L129").Value = 3999
$ws.Range("M129").Value = 3052.0001

# Row 136
$ws.Range("H136").Value = 1122666.5
$ws.Range("I136").Value = 3334666.2
$ws.Range("J136").Value = 16666.666
$ws.Range("K136").Value = 10003998.6
$ws.Range("L136").Value = 49999.99800000001
$ws.Range("M136").Value = -9998898.600000001

# Row 138
$ws.Range("H138").Value = 16674269
$ws.Range("I138").Value = 62509144
$ws.Range("J138").Value = 7041.6816
$ws.Range("K138").Value = 187527432
$ws.Range("L138").Value = 21125.0448
$ws.Range("M138").Value = -187522292

# Row 140
$ws.Range("H140").Value = 3579.8572
$ws.Range("I140").Value = 3676.5
$ws.Range("J140").Value = 3000
$ws.Range("K140").Value = 11029.5
$ws.Range("L140").Value = 9000
$ws.Range("M140").Value = -5849.5
$ws.Range("N140").Value = -19360

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2191.9524
$ws.Range("I102").Value = 2211.4211
$ws.Range("J102").Value = 2007
$ws.Range("K102").Value = 2211.4211
$ws.Range("L102").Value = 2007
$ws.Range("M102").Value = -589.4211
$ws.Range("N102").Value = -5251

# Row 126
$ws.Range("H126").Value = 11608.667
$ws.Range("I126").Value = 16151.479
$ws.Range("J126").Value = 3571.3845
$ws.Range("K126").Value = 48454.437
$ws.Range("L126").Value = 10714.1535
$ws.Range("M126").Value = -45984.437

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 6444
$ws.Range("I7").Value = 5499.5
$ws.Range("J7").Value = 8333
$ws.Range("K7").Value = 5499.5
$ws.Range("L7").Value = 8333
$ws.Range("M7").Value = -5387.5

# Row 40
$ws.Range("H40").Value = 5090.4165
$ws.Range("I40").Value = 4417.222
$ws.Range("J40").Value = 7110
$ws.Range("K40").Value = 4417.222
$ws.Range("L40").Value = 7110
$ws.Range("M40").Value = -4281.222

# Row 93
$ws.Range("H93").Value = 2549.111
$ws.Range("I93").Value = 2903.1333
$ws.Range("J93").Value = 779
$ws.Range("K93").Value = 2903.1333
$ws.Range("L93").Value = 779
$ws.Range("M93").Value = -1655.1333

# Row 122
$ws.Range("H122").Value = 4437.3
$ws.Range("I122").Value = 3911.1428
$ws.Range("J122").Value = 5665
$ws.Range("K122").Value = 11733.4284
$ws.Range("L122").Value = 16995
$ws.Range("M122").Value = -9283.428400000001
$ws.Range("N122").Value = -21895

# Row 126
$ws.Range("H126").Value = 6444
$ws.Range("I126").Value = 5499.5
$ws.Range("J126").Value = 8333
$ws.Range("K126").Value = 16498.5
$ws.Range("L126").Value = 24999
$ws.Range("M126").Value = -14028.5

# Row 132
$ws.Range("H132").Value = 3689.7585
$ws.Range("I132").Value = 2912.5
$ws.Range("J132").Value = 7420.6
$ws.Range("K132").Value = 8737.5
$ws.Range("L132").Value = 22261.8
$ws.Range("M132").Value = -6207.5

# Row 136
$ws.Range("H136").Value = 5834.579
$ws.Range("I136").Value = 5113.3335
$ws.Range("J136").Value = 7071
$ws.Range("K136").Value = 15340.0005
$ws.Range("L136").Value = 21213
$ws.Range("M136").Value = -12790.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 28
$ws.Range("H28").Value = 18233.334
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 18233.334
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 18233.334
$ws.Range("N28").Value = -18929.334

# Row 112
$ws.Range("H112").Value = 23999
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 23999
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 23999
$ws.Range("N112").Value = -26953

# Row 113
$ws.Range("H113").Value = 609.65216
$ws.Range("I113").Value = 509.35715
$ws.Range("J113").Value = 765.6667
$ws.Range("K113").Value = 1528.07145
$ws.Range("L113").Value = 2297.0001
$ws.Range("M113").Value = 641.9285500000001
$ws.Range("N113").Value = -6637.0001

# Row 122
$ws.Range("H122").Value = 3379.5908
$ws.Range("I122").Value = 3255.647
$ws.Range("J122").Value = 3801
$ws.Range("K122").Value = 9766.940999999999
$ws.Range("L122").Value = 11403
$ws.Range("M122").Value = -7316.940999999999
$ws.Range("N122").Value = -16303

# Row 132
$ws.Range("H132").Value = 3975.8845
$ws.Range("I132").Value = 4023.8635
$ws.Range("J132").Value = 3712
$ws.Range("K132").Value = 12071.5905
$ws.Range("L132").Value = 11136
$ws.Range("M132").Value = -9541.5905
